$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New work-log entry in row 23: date, hours, and description.
# Date 44047 = 2020-08-04 (serial number, matches column A's existing date formatting).
$ws.Range("A23").Value = 44047
$ws.Range("B23").Value = 6
$ws.Range("C23").Value = "User malleja toimintaan - thirdparty login ja järkevän tavan etsintää, toteuttamaan käyttäjänhallintaa"

# Row 23 grows to fit the wrapped description text, same as the other multi-line rows.
$ws.Rows.Item(23).RowHeight = 45

# Leave the cursor/selection where the author left it after the edit.
$ws.Activate() | Out-Null
$ws.Range("K20").Select() | Out-Null
